# Scheduled-runner update: refresh market-board price snapshots (and the
# derived Leve profit figures) across the per-job Sheets.
#
# Columns (all sheets share the same Table_<JOB> layout):
#   H currentAveragePrice   I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ           L LevePriceHQ
#   M LeveProfitNQ          N LeveProfitHQ

function Set-LeveRow {
    param(
        [string]$SheetName,
        [int]$Row,
        $H, $I, $J, $K, $L, $M, $N,
        [bool]$ClearM,
        [bool]$ClearN
    )

    $ws = $wb.Worksheets.Item($SheetName)

    if ($null -ne $H) { $ws.Cells.Item($Row, 8).Value  = $H }
    if ($null -ne $I) { $ws.Cells.Item($Row, 9).Value  = $I }
    if ($null -ne $J) { $ws.Cells.Item($Row, 10).Value = $J }
    if ($null -ne $K) { $ws.Cells.Item($Row, 11).Value = $K }
    if ($null -ne $L) { $ws.Cells.Item($Row, 12).Value = $L }
    if ($null -ne $M) { $ws.Cells.Item($Row, 13).Value = $M }
    if ($null -ne $N) { $ws.Cells.Item($Row, 14).Value = $N }

    if ($ClearM) { $ws.Cells.Item($Row, 13).ClearContents() }
    if ($ClearN) { $ws.Cells.Item($Row, 14).ClearContents() }
}

$wb = $excel.ActiveWorkbook

Set-LeveRow "ALC" 9 136 20 $null 20 $null 149 $null $false $false
Set-LeveRow "ALC" 19 4460.4 2150.5 $null 2150.5 $null -1975.5 $null $false $false
Set-LeveRow "ALC" 34 1521.5 1521.5 0 1521.5 0 -1318.5 $null $false $true
Set-LeveRow "ALC" 36 1521.5 1521.5 0 1521.5 0 -806.5 $null $false $true
Set-LeveRow "ALC" 38 1588.5555 137.85715 6666 413.57145 19998 -41.57144999999997 -20742 $false $false
Set-LeveRow "ALC" 40 2200 2200 2200 2200 2200 -2025 -2550 $false $false
Set-LeveRow "ALC" 51 9909.091 9000 10000 9000 10000 -8516 -10968 $false $false
Set-LeveRow "ALC" 53 1111.3334 583.6 3750 583.6 3750 53.39999999999998 -5024 $false $false
Set-LeveRow "ALC" 127 5000 5000 $null 15000 $null -10040 $null $false $false
Set-LeveRow "ALC" 138 2299.6667 697 2500 2091 7500 3049 -17780 $false $false
Set-LeveRow "ALC" 141 1104.25 1103.5 $null 3310.5 $null 1869.5 $null $false $false
Set-LeveRow "ARM" 5 280 280 $null 280 $null -168 $null $false $false
Set-LeveRow "ARM" 15 30000 $null 30000 $null 30000 $null -30700 $false $false
Set-LeveRow "ARM" 74 1492.75 1191.3 3000 1191.3 3000 -317.3 -4748 $false $false
Set-LeveRow "ARM" 77 1492.75 1191.3 3000 5956.5 15000 -1588.5 -23736 $false $false
Set-LeveRow "BSM" 4 280 280 $null 280 $null -165 $null $false $false
Set-LeveRow "CRP" 7 26.222221 18.285715 $null 18.285715 $null 94.714285 $null $false $false
Set-LeveRow "CRP" 22 487.5 475 500 475 500 -125 -1200 $false $false
Set-LeveRow "CRP" 31 1976.45 1449 $null 1449 $null -1154 $null $false $false
Set-LeveRow "CRP" 34 1976.45 1449 $null 1449 $null -1247 $null $false $false
Set-LeveRow "CUL" 4 679.9231 102.71429 $null 308.14287 $null -196.14287 $null $false $false
Set-LeveRow "CUL" 23 421 $null 592.5 $null 1777.5 $null -2247.5 $false $false
Set-LeveRow "CUL" 119 12499 12499 $null 37497 $null -32659 $null $false $false
Set-LeveRow "CUL" 124 8266.666999999999 4900 15000 14700 45000 -9790 -54820 $false $false
Set-LeveRow "CUL" 126 9999.5 4999 15000 14997 45000 -10057 -54880 $false $false
Set-LeveRow "CUL" 131 4166.6665 $null 4166.6665 $null 12499.9995 $null -22579.9995 $false $false
Set-LeveRow "CUL" 139 5005.9 5005.9 $null 15017.7 $null -9877.699999999999 $null $false $false
Set-LeveRow "GSM" 44 0 0 0 0 0 $null $null $true $true
Set-LeveRow "GSM" 47 0 0 $null 0 $null $null $null $true $false
Set-LeveRow "LTW" 2 3340 3340 $null 3340 $null -3228 $null $false $false
Set-LeveRow "LTW" 22 916.3333 749.5 1250 749.5 1250 -454.5 -1840 $false $false
Set-LeveRow "LTW" 27 916.3333 749.5 1250 749.5 1250 -642.5 -1464 $false $false
Set-LeveRow "WVR" 14 6601.8 $null 6002.25 $null 6002.25 $null -6338.25 $false $false
Set-LeveRow "WVR" 62 0 0 $null 0 $null $null $null $true $false
Set-LeveRow "WVR" 65 0 0 $null 0 $null $null $null $true $false
Set-LeveRow "WVR" 126 1385.9 1385.9 $null 4157.700000000001 $null -1687.700000000001 $null $false $false
